$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header texts (row 1 / row 2 labels) - bilingual labels + license note
$ws.Range("E1").Value = "Nivel Level"
$ws.Range("F1").Value = "#"
$ws.Range("A2").Value = "Pregunta / Question"
$ws.Range("B2").Value = "Respuesta / Answer 1"
$ws.Range("C2").Value = "Respuesta / Answer 2"
$ws.Range("D2").Value = "Respuesta / Answer 3"

# Enable wrap text on the "Nivel Level" column header/body cells
$ws.Range("E1").WrapText = $true
$ws.Range("E2").WrapText = $true

# Move the active selection
[void]$ws.Range("E4").Select()

Write-Host "done"
